$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.822.34'
$ws.Range("E2").Value = '  +4.43%  '

$ws.Range("D3").Value = '2.273.68'
$ws.Range("E3").Value = '  +1.92%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '''231.31'
$ws.Range("E5").Value = '  -0.16%  '

$ws.Range("E6").Value = '  +0.51%  '

$ws.Range("D7").Value = '''61.61'
$ws.Range("E7").Value = '  +0.82%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("E9").Value = '  +5.69%  '

$ws.Range("D10").Value = '''0.0948'
$ws.Range("E10").Value = '  +6.17%  '

$ws.Range("D11").Value = '''57.80'
$ws.Range("E11").Value = '  -2.00%  '

$ws.Range("E12").Value = '  +0.96%  '

$ws.Range("D13").Value = '2.610.78'
$ws.Range("E13").Value = '  +1.90%  '

$ws.Range("D14").Value = '''15.77'
$ws.Range("E14").Value = '  +0.76%  '

$ws.Range("D15").Value = '''23.70'
$ws.Range("E15").Value = '  +9.16%  '

$ws.Range("E16").Value = '  +4.41%  '

$ws.Range("E17").Value = '  +1.40%  '

$ws.Range("D18").Value = '2.271.54'
$ws.Range("E18").Value = '  +1.05%  '

$ws.Range("D19").Value = '43.721.44'
$ws.Range("E19").Value = '  +4.58%  '

$ws.Range("D20").Value = '0.0₃0940'
$ws.Range("E20").Value = '  +5.52%  '

$ws.Range("D21").Value = '''73.10'
$ws.Range("E21").Value = '  +1.08%  '

$ws.Range("E22").Value = '  +3.55%  '

$ws.Range("E23").Value = '  +0.42%  '

$ws.Range("E24").Value = '  -0.08%  '

$ws.Range("D25").Value = '''2.55'
$ws.Range("E25").Value = '  +7.02%  '

$ws.Range("D26").Value = '''2.49'
$ws.Range("E26").Value = '  +7.59%  '

$ws.Range("D27").Value = '''9.86'
$ws.Range("E27").Value = '  +2.31%  '

$ws.Range("D28").Value = '''170.91'
$ws.Range("E28").Value = '  +2.20%  '

$ws.Range("E29").Value = '  -1.80%  '

$ws.Range("D30").Value = '''20.61'
$ws.Range("E30").Value = '  +3.24%  '

$ws.Range("E31").Value = '  +3.98%  '

$ws.Range("E32").Value = '  +1.04%  '

$ws.Range("E33").Value = '  +0.35%  '

$ws.Range("D34").Value = '''5.12'
$ws.Range("E34").Value = '  +3.44%  '

$ws.Range("D35").Value = '''4.79'
$ws.Range("E35").Value = '  +3.33%  '

$ws.Range("D36").Value = '''0.0664'
$ws.Range("E36").Value = '  +5.10%  '

$ws.Range("E37").Value = '  -2.20%  '

$ws.Range("E38").Value = '  +2.28%  '

$ws.Range("E39").Value = '  -1.75%  '

$ws.Range("E40").Value = '  +4.06%  '

$ws.Range("E41").Value = '  +0.19%  '

$ws.Range("D42").Value = '''8.79'
$ws.Range("E42").Value = '  +2.51%  '

$ws.Range("D43").Value = '''0.000228'
$ws.Range("E43").Value = '  -9.81%  '

$ws.Range("D44").Value = '''0.0987'
$ws.Range("E44").Value = '  +0.87%  '

$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = '''1.22'
$ws.Range("E45").Value = '  +0.40%  '

$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D46").Value = '''4.48'
$ws.Range("E46").Value = '  -7.18%  '

$ws.Range("D47").Value = '''98.11'
$ws.Range("E47").Value = '  -0.86%  '

$ws.Range("D48").Value = '1.474.26'
$ws.Range("E48").Value = '  +0.00%  '

$ws.Range("D49").Value = '''16.69'
$ws.Range("E49").Value = '  +0.96%  '

$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").Value = '''1.09'
$ws.Range("E50").Value = '  +0.68%  '

$ws.Range("B51").Value = 'Celestia'
$ws.Range("C51").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D51").Value = '''9.82'
$ws.Range("E51").Value = '  +11.76%  '

